$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.304610401831553
$ws.Range("C2").Value = 0.3170129671150619
$ws.Range("E2").Value = 0.06399016100356292
$ws.Range("F2").Value = 2.994349281435461
$ws.Range("G2").Value = 0.002564756800084399
$ws.Range("I2").Value = 1.817068789092097
$ws.Range("J2").Value = 0.1265883087075839
$ws.Range("L2").Value = 0.3426189057107791
$ws.Range("M2").Value = 0.4901941093412248
$ws.Range("N2").Value = 2.266691422384085
$ws.Range("B3").Value = 2.205259143556304
$ws.Range("C3").Value = 0.2830076992580075
$ws.Range("E3").Value = 0.06278229866926743
$ws.Range("F3").Value = 2.981554429324746
$ws.Range("G3").Value = 0.002570210412247775
$ws.Range("I3").Value = 1.81538913966844
$ws.Range("J3").Value = 0.1278333547390629
$ws.Range("L3").Value = 0.3408034560854105
$ws.Range("M3").Value = 0.4762086495737563
$ws.Range("N3").Value = 2.288435265045486
$ws.Range("B4").Value = 2.145572075245241
$ws.Range("C4").Value = 0.2622657990619359
$ws.Range("E4").Value = 0.06202608509202889
$ws.Range("F4").Value = 2.975213051974677
$ws.Range("G4").Value = 0.002573735499988531
$ws.Range("I4").Value = 1.815296149735445
$ws.Range("J4").Value = 0.1286335460213266
$ws.Range("L4").Value = 0.3398117885027361
$ws.Range("M4").Value = 0.4678830974637691
$ws.Range("N4").Value = 2.302519670441797
$ws.Range("B5").Value = 2.121579046535714
$ws.Range("C5").Value = 0.2538470535120609
$ws.Range("E5").Value = 0.06171419247133159
$ws.Range("F5").Value = 2.973009208376979
$ws.Range("G5").Value = 0.002575216547592669
$ws.Range("I5").Value = 1.815493626592328
$ws.Range("J5").Value = 0.1289686031494703
$ws.Range("L5").Value = 0.3394387070239802
$ws.Range("M5").Value = 0.464556130791614
$ws.Range("N5").Value = 2.308443356777055
$ws.Range("B6").Value = 2.117614928995636
$ws.Range("C6").Value = 0.2524511416374366
$ws.Range("E6").Value = 0.06166217585958478
$ws.Range("F6").Value = 2.972666215825626
$ws.Range("G6").Value = 0.002575465169313963
$ws.Range("I6").Value = 1.815540616565038
$ws.Range("J6").Value = 0.1290247807325695
$ws.Range("L6").Value = 0.3393786341542437
$ws.Range("M6").Value = 0.4640076622140228
$ws.Range("N6").Value = 2.309438095729078
$ws.Range("B7").Value = 2.145247161933753
$ws.Range("C7").Value = 0.2621521253479386
$ws.Range("E7").Value = 0.06202189398097957
$ws.Range("F7").Value = 2.97518179103767
$ws.Range("G7").Value = 0.002573755293513035
$ws.Range("I7").Value = 1.815297860747712
$ws.Range("J7").Value = 0.1286380284051667
$ws.Range("L7").Value = 0.3398066312232686
$ws.Range("M7").Value = 0.4678379625689928
$ws.Range("N7").Value = 2.30259881402938
$ws.Range("B8").Value = 2.270080673307234
$ws.Range("C8").Value = 0.3052589120672735
$ws.Range("E8").Value = 0.06357667821630919
$ws.Range("F8").Value = 2.989622881202678
$ws.Range("G8").Value = 0.00256660065438727
$ws.Range("I8").Value = 1.816294542646929
$ws.Range("J8").Value = 0.1270101784040607
$ws.Range("L8").Value = 0.3419674602820635
$ws.Range("M8").Value = 0.4853176183417887
$ws.Range("N8").Value = 2.274036266353122
$ws.Range("B9").Value = 2.525362490340513
$ws.Range("C9").Value = 0.3909234369369301
$ws.Range("E9").Value = 0.06651284161980442
$ws.Range("F9").Value = 3.029991600802603
$ws.Range("G9").Value = 0.002553964350996358
$ws.Range("I9").Value = 1.825722804749546
$ws.Range("J9").Value = 0.12410192425516
$ws.Range("L9").Value = 0.3471773674802279
$ws.Range("M9").Value = 0.5216738357142674
$ws.Range("N9").Value = 2.223858370773499
$ws.Range("B10").Value = 2.719402045781976
$ws.Range("C10").Value = 0.4546152195821946
$ws.Range("E10").Value = 0.06860549336886379
$ws.Range("F10").Value = 3.067047696564899
$ws.Range("G10").Value = 0.00254552055276692
$ws.Range("I10").Value = 1.837249282927118
$ws.Range("J10").Value = 0.1221389219146864
$ws.Range("L10").Value = 0.3515939165796311
$ws.Range("M10").Value = 0.5496601416765188
$ws.Range("N10").Value = 2.190564775775371
$ws.Range("B11").Value = 2.809105554591156
$ws.Range("C11").Value = 0.483768257753411
$ws.Range("E11").Value = 0.06954439880551
$ws.Range("F11").Value = 3.085523511385844
$ws.Range("G11").Value = 0.002541859602316595
$ws.Range("I11").Value = 1.843501453384619
$ws.Range("J11").Value = 0.1212837511085691
$ws.Range("L11").Value = 0.3537301913441127
$ws.Range("M11").Value = 0.5626708476773672
$ws.Range("N11").Value = 2.176197835187637
$ws.Range("B12").Value = 2.843281479750374
$ws.Range("C12").Value = 0.4948345778610701
$ws.Range("E12").Value = 0.06989813357666463
$ws.Range("F12").Value = 3.092753443112684
$ws.Range("G12").Value = 0.002540499046796108
$ws.Range("I12").Value = 1.846014787410056
$ws.Range("J12").Value = 0.1209653719776327
$ws.Range("L12").Value = 0.3545573498993519
$ws.Range("M12").Value = 0.5676379689793905
$ws.Range("N12").Value = 2.170869738951993
$ws.Range("B13").Value = 2.835911854022925
$ws.Range("C13").Value = 0.4924500458374155
$ws.Range("E13").Value = 0.06982202979144247
$ws.Range("F13").Value = 3.091185948871015
$ws.Range("G13").Value = 0.002540790923150236
$ws.Range("I13").Value = 1.845467001077182
$ws.Range("J13").Value = 0.1210336976815753
$ws.Range("L13").Value = 0.3543783985537914
$ws.Range("M13").Value = 0.5665664196815499
$ws.Range("N13").Value = 2.172012236329309
$ws.Range("B14").Value = 2.811913071276763
$ws.Range("C14").Value = 0.4846781502595263
$ws.Range("E14").Value = 0.06957353662436816
$ws.Range("F14").Value = 3.086113637296791
$ws.Range("G14").Value = 0.002541747153214194
$ws.Range("I14").Value = 1.843705301218193
$ws.Range("J14").Value = 0.1212574484205788
$ws.Range("L14").Value = 0.353797877969896
$ws.Range("M14").Value = 0.5630786888108261
$ws.Range("N14").Value = 2.175757235434602
$ws.Range("B15").Value = 2.797240130258672
$ws.Range("C15").Value = 0.4799211437995154
$ws.Range("E15").Value = 0.06942109393335194
$ws.Range("F15").Value = 3.083037136797913
$ws.Range("G15").Value = 0.002542336223064874
$ws.Range("I15").Value = 1.842645214676068
$ws.Range("J15").Value = 0.1213952132406675
$ws.Range("L15").Value = 0.3534446596205072
$ws.Range("M15").Value = 0.5609475959833503
$ws.Range("N15").Value = 2.178065799338142
$ws.Range("B16").Value = 2.713568645870339
$ws.Range("C16").Value = 0.4527136860353949
$ws.Range("E16").Value = 0.06854387688836105
$ws.Range("F16").Value = 3.065872897709738
$ws.Range("G16").Value = 0.002545763417838227
$ws.Range("I16").Value = 1.836861048183238
$ws.Range("J16").Value = 0.122195571724907
$ws.Range("L16").Value = 0.3514568572861521
$ws.Range("M16").Value = 0.5488154912439072
$ws.Range("N16").Value = 2.191519388849578
$ws.Range("B17").Value = 2.662606835241377
$ws.Range("C17").Value = 0.4360693658341006
$ws.Range("E17").Value = 0.06800243817861507
$ws.Range("F17").Value = 3.055758343959226
$ws.Range("G17").Value = 0.002547911937104979
$ws.Range("I17").Value = 1.833571515146247
$ws.Range("J17").Value = 0.1226962625106953
$ws.Range("L17").Value = 0.3502699114921484
$ws.Range("M17").Value = 0.5414444785670582
$ws.Range("N17").Value = 2.199972405824127
$ws.Range("B18").Value = 2.633429859910962
$ws.Range("C18").Value = 0.4265128170260368
$ws.Range("E18").Value = 0.06768978425869321
$ws.Range("F18").Value = 3.050093044019292
$ws.Range("G18").Value = 0.002549164676052529
$ws.Range("I18").Value = 1.831774364296152
$ws.Range("J18").Value = 0.1229878041954944
$ws.Range("L18").Value = 0.34959918619424
$ws.Range("M18").Value = 0.5372311809090462
$ws.Range("N18").Value = 2.204907607280177
$ws.Range("B19").Value = 2.623574179836851
$ws.Range("C19").Value = 0.4232799957910629
$ws.Range("E19").Value = 0.0675837112208626
$ws.Range("F19").Value = 3.048201009699639
$ws.Range("G19").Value = 0.002549591750224323
$ws.Range("I19").Value = 1.8311821575506
$ws.Range("J19").Value = 0.1230871256624599
$ws.Range("L19").Value = 0.349374149483296
$ws.Range("M19").Value = 0.5358091488936623
$ws.Range("N19").Value = 2.206591150090198
$ws.Range("B20").Value = 2.668017837482296
$ws.Range("C20").Value = 0.4378394343615355
$ws.Range("E20").Value = 0.06806020233833188
$ws.Range("F20").Value = 3.056819285092672
$ws.Range("G20").Value = 0.00254768146852997
$ws.Range("I20").Value = 1.833911864424977
$ws.Range("J20").Value = 0.1226425947921763
$ws.Range("L20").Value = 0.3503950252245147
$ws.Range("M20").Value = 0.542226411762158
$ws.Range("N20").Value = 2.199064983750283
$ws.Range("B21").Value = 2.818956465654935
$ws.Range("C21").Value = 0.4869602125165784
$ws.Range("E21").Value = 0.06964657361877258
$ws.Range("F21").Value = 3.087597153781189
$ws.Range("G21").Value = 0.002541465587031433
$ws.Range("I21").Value = 1.844218793141167
$ws.Range("J21").Value = 0.1211915791920593
$ws.Range("L21").Value = 0.3539678977569025
$ws.Range("M21").Value = 0.5641020271817681
$ws.Range("N21").Value = 2.174654185919344
$ws.Range("B22").Value = 2.91881149943066
$ws.Range("C22").Value = 0.5192195235928807
$ws.Range("E22").Value = 0.07067285093107856
$ws.Range("F22").Value = 3.109073893119842
$ws.Range("G22").Value = 0.00253755327326366
$ws.Range("I22").Value = 1.851804984602353
$ws.Range("J22").Value = 0.1202750645525086
$ws.Range("L22").Value = 0.3564090062610461
$ws.Range("M22").Value = 0.5786336342148104
$ws.Range("N22").Value = 2.159355388549017
$ws.Range("B23").Value = 2.865406140091409
$ws.Range("C23").Value = 0.5019875403547189
$ws.Range("E23").Value = 0.07012604604534367
$ws.Range("F23").Value = 3.09748650966597
$ws.Range("G23").Value = 0.002539627658934571
$ws.Range("I23").Value = 1.847678073736645
$ws.Range("J23").Value = 0.120761308896725
$ws.Range("L23").Value = 0.3550964678745459
$ws.Range("M23").Value = 0.5708563577736498
$ws.Range("N23").Value = 2.167460568989142
$ws.Range("B24").Value = 2.665571144532692
$ws.Range("C24").Value = 0.4370391475371775
$ws.Range("E24").Value = 0.06803409144532857
$ws.Range("F24").Value = 3.056339167406833
$ws.Range("G24").Value = 0.002547785608745318
$ws.Range("I24").Value = 1.833757699639605
$ws.Range("J24").Value = 0.1226668464823666
$ws.Range("L24").Value = 0.3503384249635815
$ws.Range("M24").Value = 0.5418728237805865
$ws.Range("N24").Value = 2.199474994292615
$ws.Range("B25").Value = 2.455170164585411
$ws.Range("C25").Value = 0.3676211422929327
$ws.Range("E25").Value = 0.0657301725230397
$ws.Range("F25").Value = 3.017775530499478
$ws.Range("G25").Value = 0.002557234579861104
$ws.Range("I25").Value = 1.822367652554433
$ws.Range("J25").Value = 0.1248582266826279
$ws.Range("L25").Value = 0.3456642567984005
$ws.Range("M25").Value = 0.5116151521946861
$ws.Range("N25").Value = 2.236806708990954
